# Apply cryptos list update (prices & volume %) as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value (avoids Excel auto-converting
    # numeric-looking strings like "1.001" into a float), then restore
    # the default (no explicit) cell style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '27.730.66'
$ws.Range("E2").Value = '  -0.12%  '
# Row 3
Set-TextValue $ws.Range("D3") '1.900.58'
$ws.Range("E3").Value = '  +0.20%  '
# Row 4
$ws.Range("E4").Value = '  +0.14%  '
# Row 5
Set-TextValue $ws.Range("D5") '312.04'
$ws.Range("E5").Value = '  +0.03%  '
# Row 6
Set-TextValue $ws.Range("D6") '1.001'
$ws.Range("E6").Value = '  +0.12%  '
# Row 7
Set-TextValue $ws.Range("D7") '0.5221'
$ws.Range("E7").Value = '  +6.06%  '
# Row 8
Set-TextValue $ws.Range("D8") '0.3788'
$ws.Range("E8").Value = '  -0.26%  '
# Row 9
Set-TextValue $ws.Range("D9") '0.07248'
$ws.Range("E9").Value = '  -1.04%  '
# Row 10
Set-TextValue $ws.Range("D10") '21.23'
$ws.Range("E10").Value = '  +2.78%  '
# Row 11
Set-TextValue $ws.Range("D11") '0.9023'
$ws.Range("E11").Value = '  -1.02%  '
# Row 12
Set-TextValue $ws.Range("D12") '0.07646'
$ws.Range("E12").Value = '  +0.28%  '
# Row 13
Set-TextValue $ws.Range("D13") '1.891.88'
$ws.Range("E13").Value = '  -0.18%  '
# Row 14
Set-TextValue $ws.Range("D14") '5.448'
$ws.Range("E14").Value = '  -0.28%  '
# Row 15
Set-TextValue $ws.Range("D15") '92.24'
$ws.Range("E15").Value = '  +1.16%  '
# Row 16
$ws.Range("E16").Value = '  +0.13%  '
# Row 17
Set-TextValue $ws.Range("D17") '0.000008704'
$ws.Range("E17").Value = '  -0.38%  '
# Row 18
Set-TextValue $ws.Range("D18") '1.001'
$ws.Range("E18").Value = '  +0.12%  '
# Row 19
Set-TextValue $ws.Range("D19") '27.760.31'
$ws.Range("E19").Value = '  -0.06%  '
# Row 20
Set-TextValue $ws.Range("D20") '14.45'
$ws.Range("E20").Value = '  -0.20%  '
# Row 21
Set-TextValue $ws.Range("D21") '5.141'
$ws.Range("E21").Value = '  +0.44%  '
# Row 22
Set-TextValue $ws.Range("D22") '2.132.05'
$ws.Range("E22").Value = '  +0.60%  '
# Row 23
Set-TextValue $ws.Range("D23") '10.85'
$ws.Range("E23").Value = '  +0.76%  '
# Row 24
Set-TextValue $ws.Range("D24") '6.621'
$ws.Range("E24").Value = '  -0.50%  '
# Row 25
Set-TextValue $ws.Range("D25") '153.18'
$ws.Range("E25").Value = '  -0.51%  '
# Row 26
Set-TextValue $ws.Range("D26") '1.864'
$ws.Range("E26").Value = '  +0.69%  '
# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D27") '18.30'
$ws.Range("E27").Value = '  -0.57%  '
# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D28") '2.158'
$ws.Range("E28").Value = '  -0.81%  '
# Row 29
Set-TextValue $ws.Range("D29") '114.47'
$ws.Range("E29").Value = '  -0.82%  '
# Row 30
Set-TextValue $ws.Range("D30") '4.845'
$ws.Range("E30").Value = '  -0.86%  '
# Row 31
Set-TextValue $ws.Range("D31") '0.09094'
$ws.Range("E31").Value = '  +1.76%  '
# Row 32
Set-TextValue $ws.Range("D32") '3.192'
$ws.Range("E32").Value = '  -1.14%  '
# Row 33
Set-TextValue $ws.Range("D33") '4.835'
$ws.Range("E33").Value = '  +4.12%  '
# Row 34
Set-TextValue $ws.Range("D34") '1.226'
$ws.Range("E34").Value = '  -0.50%  '
# Row 35
Set-TextValue $ws.Range("D35") '0.7791'
$ws.Range("E35").Value = '  +1.57%  '
# Row 36
Set-TextValue $ws.Range("D36") '0.02092'
$ws.Range("E36").Value = '  +2.41%  '
# Row 37
Set-TextValue $ws.Range("D37") '2.578'
$ws.Range("E37").Value = '  +0.58%  '
# Row 38
Set-TextValue $ws.Range("D38") '3.069'
$ws.Range("E38").Value = '  +2.67%  '
# Row 39
$ws.Range("E39").Value = '  -0.62%  '
# Row 40
Set-TextValue $ws.Range("D40") '0.5551'
$ws.Range("E40").Value = '  +0.86%  '
# Row 41
Set-TextValue $ws.Range("D41") '0.05291'
$ws.Range("E41").Value = '  +0.12%  '
# Row 42
Set-TextValue $ws.Range("D42") '6.725'
$ws.Range("E42").Value = '  -2.53%  '
# Row 43
Set-TextValue $ws.Range("D43") '116.22'
$ws.Range("E43").Value = '  +3.01%  '
# Row 44
Set-TextValue $ws.Range("D44") '8.505'
$ws.Range("E44").Value = '  -0.83%  '
# Row 45
$ws.Range("E45").Value = '  -0.39%  '
# Row 46
Set-TextValue $ws.Range("D46") '0.4813'
$ws.Range("E46").Value = '  +0.24%  '
# Row 47
Set-TextValue $ws.Range("D47") '10.42'
$ws.Range("E47").Value = '  -2.37%  '
# Row 48
Set-TextValue $ws.Range("D48") '1.000'
$ws.Range("E48").Value = '  +0.10%  '
# Row 49
Set-TextValue $ws.Range("D49") '1.612'
$ws.Range("E49").Value = '  -1.41%  '
# Row 50
$ws.Range("E50").Value = '  -1.16%  '
# Row 51
$ws.Range("E51").Value = '  -0.87%  '

Write-Output "Applied cryptos list update"
